$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer-rules")
$ws.Activate()

# F11: 30 -> 35
$ws.Range("F11").Value = 35

# sheetView topLeftCell: B1 -> B2 (scroll so B2 is the top-left visible cell)
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2

# keep the existing selection (activeCell G11, sqref G11)
$ws.Range("G11").Select()
